# Fruta / hortaliza, semanal
# Insert 3 new price records (Durazno, variety "Springcrest") above the
# existing row 99, pushing the former rows 99:121 down to 102:124.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the current row 99 (this shifts the
# old rows 99:121 down to 102:124, preserving all of their original data).
$ws.Rows("99:101").Insert()

# Fill in the three new rows with the new "Springcrest" records.
$newRows = @(
    @{ Row=99;  D=44559; K="Springcrest"; L="Especial"; M=50; N=15000; O=15000; P=15000; R="Región Metropolitana"; S=938 },
    @{ Row=100; D=44559; K="Springcrest"; L="Primera";  M=50; N=13000; O=13000; P=13000; R="Región Metropolitana"; S=812 },
    @{ Row=101; D=44559; K="Springcrest"; L="Segunda";  M=50; N=11000; O=11000; P=11000; R="Región Metropolitana"; S=688 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 11
    $ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($row, 3).Value = "Bíobío"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 8
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100103
    $ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($row, 9).Value = 100103004
    $ws.Cells.Item($row, 10).Value = "Durazno"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/caja 16 kilos empedrada"
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 16
}
